# Regenerate handback-status report values (new GUID-named files, new hash,
# new timestamps) while keeping hyperlink target addresses (rIds) the same.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "1d249880-a870-4b19-b3d8-9119531c964d"
$oldGuid2 = "6ec87331-1908-4f20-84f6-7c7d4a24c069"
$newGuid1 = "c316ba12-b630-4122-9583-8507e8eee9ef"
$newGuid2 = "ffff6b803c36-591f-47e8-956a-8ca7516ddc5a"

$newHash = "fb8b5184faab6789b8c930836fb6e0ce9677c175"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$addrSrc1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66aea97287dcd3b4beb48c8f9e7682292ac88436/e2e/$oldGuid1.md"
$addrSrc2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66aea97287dcd3b4beb48c8f9e7682292ac88436/e2e/$oldGuid2.md"

$ws1.Range("A2").Value = "$newGuid1.md"
$ws1.Range("G2").Value = "2016-08-17 03:04:39"

$ws1.Range("A3").Value = "$newGuid2.md"
$ws1.Range("G3").Value = "2016-08-17 03:04:39"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $addrSrc1, "", "", "e2e\$newGuid1.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), $addrSrc2, "", "", "e2e\$newGuid2.md")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$addrSrc1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66aea97287dcd3b4beb48c8f9e7682292ac88436/e2e/$oldGuid1.md"
$addrSrc2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66aea97287dcd3b4beb48c8f9e7682292ac88436/e2e/$oldGuid2.md"
$addrTgt1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d9e1872ffe31f74a958f51202b6442933e8da51d/e2e/$oldGuid1.md"
$addrTgt2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d9e1872ffe31f74a958f51202b6442933e8da51d/e2e/$oldGuid2.md"

$ws2.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-17 03:04:34"
$ws2.Range("J2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-17 03:04:50"

$ws2.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$ws2.Range("H3").Value = "2016-08-17 03:04:34"
$ws2.Range("J3").Value = "$newGuid1.$newHash.de-de.xlf"
$ws2.Range("K3").Value = "2016-08-17 03:04:50"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $addrSrc1, "", "", "$newGuid1.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), $addrTgt1, "", "", "$newGuid1.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $addrSrc2, "", "", "$newGuid2.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), $addrTgt2, "", "", "$newGuid2.md")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$addrSrc1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66aea97287dcd3b4beb48c8f9e7682292ac88436/e2e/$oldGuid1.md"
$addrSrc2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66aea97287dcd3b4beb48c8f9e7682292ac88436/e2e/$oldGuid2.md"
$addrTgt1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a2470cc69c100a0058635db2f2115285b2ec4ef1/e2e/$oldGuid1.md"
$addrTgt2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a2470cc69c100a0058635db2f2115285b2ec4ef1/e2e/$oldGuid2.md"

$ws3.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-17 03:04:39"
$ws3.Range("J2").Value = "$newGuid1.$newHash.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-17 03:04:58"

$ws3.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-17 03:04:39"
$ws3.Range("J3").Value = "$newGuid1.$newHash.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-17 03:04:58"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $addrSrc1, "", "", "$newGuid1.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), $addrTgt1, "", "", "$newGuid1.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $addrSrc2, "", "", "$newGuid2.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), $addrTgt2, "", "", "$newGuid2.md")
